# New changes as of 04/15
# The study's phs accession number was updated from phs001713 to phs002050
# throughout the workbook (TSV/Web data file names + the embedded SQL
# queries), and the active selection moved to C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldAccession = "phs001713"
$newAccession = "phs002050"

# Every cell whose text references the old phs accession number.
$cellsToUpdate = @("B2", "C2", "D2", "E2", "B3", "B4")

foreach ($cellRef in $cellsToUpdate) {
    $cell = $ws.Range($cellRef)
    $oldValue = $cell.Value2
    if ($oldValue -ne $null -and $oldValue -like "*$oldAccession*") {
        $cell.Value = $oldValue -replace $oldAccession, $newAccession
    }
}

# Rewriting the long wrapped query in B2 makes the host recompute row 2's
# auto-fit height; restore it to the original (Excel's max) row height so
# the row keeps its prior, non-autofit size.
$ws.Rows(2).RowHeight = 409.5

# Move the live selection/active cell to C13, matching the saved view state.
$ws.Activate()
$ws.Range("C13").Select()
